# Update TPM-derived NATMI metrics for Bmp2-Bmpr1b ligand/receptor pairs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7170026666666667
$ws.Range("H2").Value = 2.151008
$ws.Range("I2").Value = 0.02953485643833859
$ws.Range("J2").Value = 0.02953485643833859
$ws.Range("Q2").Value = 0.9268953982968888
$ws.Range("R2").Value = 8.342058584671999
$ws.Range("S2").Value = 0.02665310832609218
$ws.Range("T2").Value = 0.02665310832609219

# Row 3 (ECs -> FAPs)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7170026666666667
$ws.Range("H3").Value = 2.151008
$ws.Range("I3").Value = 0.02953485643833859
$ws.Range("J3").Value = 0.02953485643833859
$ws.Range("Q3").Value = 0.1002164187235556
$ws.Range("R3").Value = 0.9019477685120001
$ws.Range("S3").Value = 0.002881748112246405
$ws.Range("T3").Value = 0.002881748112246406

# Row 4 (FAPs -> ECs)
$ws.Range("I4").Value = 0.4970672037825566
$ws.Range("J4").Value = 0.4970672037825566
$ws.Range("S4").Value = 0.4485678153006616
$ws.Range("T4").Value = 0.4485678153006616

# Row 5 (FAPs -> FAPs)
$ws.Range("I5").Value = 0.4970672037825566
$ws.Range("J5").Value = 0.4970672037825566
$ws.Range("Q5").Value = 1.686627295176444
$ws.Range("S5").Value = 0.04849938848189503
$ws.Range("T5").Value = 0.04849938848189503

# Row 6 (MuSCs -> ECs)
$ws.Range("I6").Value = 0.4733979397791048
$ws.Range("J6").Value = 0.4733979397791048
$ws.Range("S6").Value = 0.4272079871667429
$ws.Range("T6").Value = 0.4272079871667429

# Row 7 (MuSCs -> FAPs)
$ws.Range("I7").Value = 0.4733979397791048
$ws.Range("J7").Value = 0.4733979397791048
$ws.Range("S7").Value = 0.04618995261236195
$ws.Range("T7").Value = 0.04618995261236195
